$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 114-117 use the same sparse column layout as row 113.
# Copy each populated cell from row 113 down to the new rows first so they
# inherit identical per-cell formatting (e.g. the date style on column A),
# then overwrite with the real values for each new row.
$cols = @("A","B","C","I","K","N","O","Q","U","Z")
foreach ($r in 114..117) {
  foreach ($c in $cols) {
    $ws.Range("$c" + "113").Copy($ws.Range("$c$r"))
  }
}

# Now set the actual values for the populated columns in each new row

# Row 114
$ws.Range("A114").Value = 45544
$ws.Range("B114").Value = 580.27742802
$ws.Range("C114").Value = 163.599117825
$ws.Range("I114").Value = 241.853053948
$ws.Range("K114").Value = 300.154098790958
$ws.Range("N114").Value = 37.24026124448
$ws.Range("O114").Value = 0.020810244
$ws.Range("Q114").Value = 0.0000017712
$ws.Range("U114").Value = 247.6705791037456
$ws.Range("Z114").Value = 201.904638488816

# Row 115
$ws.Range("A115").Value = 45545
$ws.Range("B115").Value = 586.3199754318999
$ws.Range("C115").Value = 165.611258702
$ws.Range("I115").Value = 242.837797653
$ws.Range("K115").Value = 297.241866267113
$ws.Range("N115").Value = 36.04121563392
$ws.Range("O115").Value = 0.02074204
$ws.Range("Q115").Value = 0.0000017904
$ws.Range("U115").Value = 257.0093974273889
$ws.Range("Z115").Value = 219.007442779604

# Row 116
$ws.Range("A116").Value = 45546
$ws.Range("B116").Value = 583.28857978
$ws.Range("C116").Value = 162.2851939925
$ws.Range("I116").Value = 237.090475302
$ws.Range("K116").Value = 297.824312771882
$ws.Range("N116").Value = 33.59655953472
$ws.Range("O116").Value = 0.021267612
$ws.Range("Q116").Value = 0.0000017376
$ws.Range("U116").Value = 255.3463201916716
$ws.Range("Z116").Value = 222.440487436514

# Row 117
$ws.Range("A117").Value = 45547
$ws.Range("B117").Value = 591.3690462192
$ws.Range("C117").Value = 163.755817976
$ws.Range("I117").Value = 243.589783755
$ws.Range("K117").Value = 295.10622908296
$ws.Range("N117").Value = 34.94694099904
$ws.Range("O117").Value = 0.02184534
$ws.Range("Q117").Value = 0.0000017808
$ws.Range("U117").Value = 257.137326445521
$ws.Range("Z117").Value = 241.790375502734

Write-Host "done"